$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 204078
$ws.Range("E2").Value = 3761
$ws.Range("F2").Value = 3761
$ws.Range("G2").Value = 2872
$ws.Range("H2").Value = 1764
$ws.Range("I2").Value = 1845
$ws.Range("J2").Value = -81
$ws.Range("K2").Value = 93409
$ws.Range("L2").Value = 69661
$ws.Range("M2").Value = 23747
$ws.Range("N2").Value = 23785
$ws.Range("O2").Value = -38
$ws.Range("P2").Value = 5694
$ws.Range("Q2").Value = -4376
$ws.Range("R2").Value = -4274
$ws.Range("S2").Value = 9525
$ws.Range("T2").Value = 960
$ws.Range("U2").Value = -5336
$ws.Range("V2").Value = 50319
$ws.Range("W2").Value = 1.84
$ws.Range("X2").Value = 0.86
$ws.Range("Y2").Value = 8.01
$ws.Range("Z2").Value = 2
$ws.Range("AA2").Value = 293.34
$ws.Range("AB2").Value = 304.36
$ws.Range("AC2").Value = 1620
$ws.Range("AD2").Value = 19.41
$ws.Range("AE2").Value = 20887
$ws.Range("AF2").Value = 1.51
$ws.Range("AG2").Value = 500
$ws.Range("AH2").Value = 1.59
$ws.Range("AI2").Value = 30.86
$ws.Range("AJ2").Value = 113876291

# Row 3
$ws.Range("D3").Value = 175269
$ws.Range("E3").Value = 3688
$ws.Range("F3").Value = 3688
$ws.Range("G3").Value = 1335
$ws.Range("H3").Value = 1086
$ws.Range("I3").Value = 1289
$ws.Range("J3").Value = -204
$ws.Range("K3").Value = 80433
$ws.Range("L3").Value = 56221
$ws.Range("M3").Value = 24212
$ws.Range("N3").Value = 24417
$ws.Range("O3").Value = -206
$ws.Range("P3").Value = 5694
$ws.Range("Q3").Value = 12676
$ws.Range("R3").Value = -2081
$ws.Range("S3").Value = -10960
$ws.Range("T3").Value = 1029
$ws.Range("U3").Value = 11647
$ws.Range("V3").Value = 40383
$ws.Range("W3").Value = 2.1
$ws.Range("X3").Value = 0.62
$ws.Range("Y3").Value = 5.35
$ws.Range("Z3").Value = 1.25
$ws.Range("AA3").Value = 232.2
$ws.Range("AB3").Value = 317.04
$ws.Range("AC3").Value = 1132
$ws.Range("AD3").Value = 14.44
$ws.Range("AE3").Value = 21442
$ws.Range("AF3").Value = 0.76
$ws.Range("AG3").Value = 500
$ws.Range("AH3").Value = 3.06
$ws.Range("AI3").Value = 44.16
$ws.Range("AJ3").Value = 113876291

# Row 4
$ws.Range("D4").Value = 164921
$ws.Range("E4").Value = 3181
$ws.Range("F4").Value = 3181
$ws.Range("G4").Value = 1746
$ws.Range("H4").Value = 1222
$ws.Range("I4").Value = 1113
$ws.Range("J4").Value = 110
$ws.Range("K4").Value = 82888
$ws.Range("L4").Value = 57851
$ws.Range("M4").Value = 25037
$ws.Range("N4").Value = 25158
$ws.Range("O4").Value = -122
$ws.Range("P4").Value = 5694
$ws.Range("Q4").Value = 7928
$ws.Range("R4").Value = -1828
$ws.Range("S4").Value = -6615
$ws.Range("T4").Value = 956
$ws.Range("U4").Value = 6973
$ws.Range("V4").Value = 35255
$ws.Range("W4").Value = 1.93
$ws.Range("X4").Value = 0.74
$ws.Range("Y4").Value = 4.49
$ws.Range("Z4").Value = 1.5
$ws.Range("AA4").Value = 231.07
$ws.Range("AB4").Value = 326.32
$ws.Range("AC4").Value = 977
$ws.Range("AD4").Value = 27.64
$ws.Range("AE4").Value = 22093
$ws.Range("AF4").Value = 1.22
$ws.Range("AG4").Value = 500
$ws.Range("AH4").Value = 1.85
$ws.Range("AI4").Value = 51.18
$ws.Range("AJ4").Value = 113876291

# Row 5
$ws.Range("D5").Value = 225717
$ws.Range("E5").Value = 4013
$ws.Range("F5").Value = 4013
$ws.Range("G5").Value = 2480
$ws.Range("H5").Value = 1668
$ws.Range("I5").Value = 1759
$ws.Range("J5").Value = -91
$ws.Range("K5").Value = 91700
$ws.Range("L5").Value = 63175
$ws.Range("M5").Value = 28525
$ws.Range("N5").Value = 28098
$ws.Range("O5").Value = 427
$ws.Range("P5").Value = 6169
$ws.Range("Q5").Value = 3822
$ws.Range("R5").Value = -2316
$ws.Range("S5").Value = -1465
$ws.Range("T5").Value = 728
$ws.Range("U5").Value = 3095
$ws.Range("V5").Value = 40328
$ws.Range("W5").Value = 1.78
$ws.Range("X5").Value = 0.74
$ws.Range("Y5").Value = 6.6
$ws.Range("Z5").Value = 1.91
$ws.Range("AA5").Value = 221.47
$ws.Range("AB5").Value = 356.28
$ws.Range("AC5").Value = 1444
$ws.Range("AD5").Value = 12.57
$ws.Range("AE5").Value = 22774
$ws.Range("AF5").Value = 0.8
$ws.Range("AG5").Value = 500
$ws.Range("AH5").Value = 2.75
$ws.Range("AI5").Value = 35.07
$ws.Range("AJ5").Value = 123375149

# Row 6
$ws.Range("D6").Value = 251739
$ws.Range("E6").Value = 4726
$ws.Range("F6").Value = 4726
$ws.Range("G6").Value = 1761
$ws.Range("H6").Value = 1157
$ws.Range("I6").Value = 1552
$ws.Range("K6").Value = 98810
$ws.Range("L6").Value = 69915
$ws.Range("M6").Value = 28895
$ws.Range("N6").Value = 28873
$ws.Range("P6").Value = 6169
$ws.Range("Q6").Value = -1461
$ws.Range("R6").Value = -958
$ws.Range("S6").Value = 2668
$ws.Range("T6").Value = 990
$ws.Range("U6").Value = -2451
$ws.Range("V6").Value = 43894
$ws.Range("W6").Value = 1.88
$ws.Range("X6").Value = 0.46
$ws.Range("Y6").Value = 5.45
$ws.Range("Z6").Value = 1.21
$ws.Range("AA6").Value = 241.96
$ws.Range("AB6").Value = 370.62
$ws.Range("AC6").Value = 1258
$ws.Range("AD6").Value = 14.47
$ws.Range("AE6").Value = 23402
$ws.Range("AF6").Value = 0.78
$ws.Range("AG6").Value = 600
$ws.Range("AH6").Value = 3.3
$ws.Range("AI6").Value = 47.7
$ws.Range("AJ6").Value = 123375149

# Row 7
$ws.Range("D7").Value = 244315
$ws.Range("E7").Value = 6382
$ws.Range("G7").Value = 3841
$ws.Range("H7").Value = 2533
$ws.Range("I7").Value = 2579
$ws.Range("K7").Value = 98253
$ws.Range("L7").Value = 67308
$ws.Range("M7").Value = 30945
$ws.Range("N7").Value = 30956
$ws.Range("P7").Value = 6170
$ws.Range("Q7").Value = 8478
$ws.Range("R7").Value = -853
$ws.Range("S7").Value = -5592
$ws.Range("T7").Value = 1163
$ws.Range("U7").Value = 7358
$ws.Range("W7").Value = 2.61
$ws.Range("X7").Value = 1.04
$ws.Range("Y7").Value = 8.619999999999999
$ws.Range("Z7").Value = 2.57
$ws.Range("AA7").Value = 217.51
$ws.Range("AC7").Value = 2091
$ws.Range("AD7").Value = 8.470000000000001
$ws.Range("AE7").Value = 25091
$ws.Range("AF7").Value = 0.71
$ws.Range("AG7").Value = 625
$ws.Range("AH7").Value = 3.53
$ws.Range("AI7").Value = 29.89

# Row 8
$ws.Range("D8").Value = 246656
$ws.Range("E8").Value = 6104
$ws.Range("G8").Value = 4843
$ws.Range("H8").Value = 3380
$ws.Range("I8").Value = 3368
$ws.Range("K8").Value = 99147
$ws.Range("L8").Value = 65676
$ws.Range("M8").Value = 33471
$ws.Range("N8").Value = 33474
$ws.Range("P8").Value = 6170
$ws.Range("Q8").Value = 6837
$ws.Range("R8").Value = -2145
$ws.Range("S8").Value = -2810
$ws.Range("T8").Value = 1194
$ws.Range("U8").Value = 5138
$ws.Range("W8").Value = 2.48
$ws.Range("X8").Value = 1.37
$ws.Range("Y8").Value = 10.46
$ws.Range("Z8").Value = 3.43
$ws.Range("AA8").Value = 196.21
$ws.Range("AC8").Value = 2730
$ws.Range("AD8").Value = 6.17
$ws.Range("AE8").Value = 27132
$ws.Range("AF8").Value = 0.62
$ws.Range("AG8").Value = 717
$ws.Range("AH8").Value = 4.25
$ws.Range("AI8").Value = 26.25

# Row 9
$ws.Range("D9").Value = 255204
$ws.Range("E9").Value = 6157
$ws.Range("G9").Value = 5016
$ws.Range("H9").Value = 3501
$ws.Range("I9").Value = 3490
$ws.Range("K9").Value = 103882
$ws.Range("L9").Value = 67837
$ws.Range("M9").Value = 36045
$ws.Range("N9").Value = 36037
$ws.Range("P9").Value = 6170
$ws.Range("Q9").Value = 5342
$ws.Range("R9").Value = -2076
$ws.Range("S9").Value = -2081
$ws.Range("T9").Value = 1086
$ws.Range("U9").Value = 4311
$ws.Range("W9").Value = 2.41
$ws.Range("X9").Value = 1.37
$ws.Range("Y9").Value = 10.04
$ws.Range("Z9").Value = 3.45
$ws.Range("AA9").Value = 188.2
$ws.Range("AC9").Value = 2828
$ws.Range("AD9").Value = 5.96
$ws.Range("AE9").Value = 29209
$ws.Range("AF9").Value = 0.58
$ws.Range("AG9").Value = 733
$ws.Range("AH9").Value = 4.35
$ws.Range("AI9").Value = 25.93
